$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Update the date in A1 (one month later: 2024-04-24 -> 2024-05-24)
$ws.Range("A1").Value = "05/24/2024"

# Update PRECIO column (D29:D38) with the new prices
$ws.Range("D29").Value = 243.87
$ws.Range("D30").Value = 243.87
$ws.Range("D31").Value = 243.87
$ws.Range("D32").Value = 243.87
$ws.Range("D33").Value = 262.13
$ws.Range("D34").Value = 262.13
$ws.Range("D35").Value = 262.13
$ws.Range("D36").Value = 262.13
$ws.Range("D37").Value = 289.3
$ws.Range("D38").Value = 289.3
